$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '27.741.62'
$ws.Cells.Item(2, 5).Value = '  +1.31%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.877.92'
$ws.Cells.Item(3, 5).Value = '  +1.05%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.002'
$ws.Cells.Item(4, 5).Value = '  -0.03%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '331.94'
$ws.Cells.Item(5, 5).Value = '  +2.56%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  +0.09%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.4720'
$ws.Cells.Item(7, 5).Value = '  +4.21%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3951'
$ws.Cells.Item(8, 5).Value = '  +2.01%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '47.86'
$ws.Cells.Item(9, 5).Value = '  -2.17%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.08054'
$ws.Cells.Item(10, 5).Value = '  +1.69%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +1.17%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '22.13'
$ws.Cells.Item(12, 5).Value = '  +3.45%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.884.37'
$ws.Cells.Item(13, 5).Value = '  +1.60%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '5.965'
$ws.Cells.Item(14, 5).Value = '  +0.75%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '7.136'
$ws.Cells.Item(15, 5).Value = '  +0.05%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '1.005'
$ws.Cells.Item(16, 5).Value = '  +0.12%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.00001049'
$ws.Cells.Item(17, 5).Value = '  +1.48%  '

# Row 18
$ws.Cells.Item(18, 2).Value = 'Litecoin'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '87.05'
$ws.Cells.Item(18, 5).Value = '  +1.18%  '

# Row 19
$ws.Cells.Item(19, 2).Value = 'TRON'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '0.06659'
$ws.Cells.Item(19, 5).Value = '  +2.18%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '17.27'
$ws.Cells.Item(20, 5).Value = '  +1.25%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '1.003'
$ws.Cells.Item(21, 5).Value = '  +0.10%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '27.747.29'
$ws.Cells.Item(22, 5).Value = '  +1.36%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '5.527'
$ws.Cells.Item(23, 5).Value = '  -0.28%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '10.98'
$ws.Cells.Item(24, 5).Value = '  +1.15%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '2.304'
$ws.Cells.Item(25, 5).Value = '  +1.00%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '2.111.70'
$ws.Cells.Item(26, 5).Value = '  +1.67%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '158.83'
$ws.Cells.Item(27, 5).Value = '  +3.18%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '20.17'
$ws.Cells.Item(28, 5).Value = '  +1.23%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '2.108'
$ws.Cells.Item(29, 5).Value = '  +1.41%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '5.584'
$ws.Cells.Item(30, 5).Value = '  +2.59%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '122.45'
$ws.Cells.Item(31, 5).Value = '  +1.10%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '0.9763'
$ws.Cells.Item(32, 5).Value = '  +4.13%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.09537'
$ws.Cells.Item(33, 5).Value = '  +2.58%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '1.453'
$ws.Cells.Item(34, 5).Value = '  -2.08%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '3.591'
$ws.Cells.Item(35, 5).Value = '  -0.19%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '5.341'
$ws.Cells.Item(36, 5).Value = '  +1.42%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.06111'
$ws.Cells.Item(37, 5).Value = '  +1.82%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.02257'
$ws.Cells.Item(38, 5).Value = '  +0.91%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '1.231'
$ws.Cells.Item(39, 5).Value = '  -0.15%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '8.156'
$ws.Cells.Item(40, 5).Value = '  -0.53%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.6035'
$ws.Cells.Item(41, 5).Value = '  +2.01%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.1903'
$ws.Cells.Item(42, 5).Value = '  -0.05%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '10.28'
$ws.Cells.Item(43, 5).Value = '  +1.61%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '1.266'
$ws.Cells.Item(44, 5).Value = '  -1.37%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.5728'
$ws.Cells.Item(45, 5).Value = '  +1.88%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '12.17'
$ws.Cells.Item(46, 5).Value = '  +1.28%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '1.951'
$ws.Cells.Item(47, 5).Value = '  +1.32%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '3.377'
$ws.Cells.Item(48, 5).Value = '  +0.13%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.06890'
$ws.Cells.Item(49, 5).Value = '  +1.80%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '114.91'
$ws.Cells.Item(50, 5).Value = '  +5.94%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '1.071'
$ws.Cells.Item(51, 5).Value = '  +1.65%  '
